$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 269
$ws1.Range("F5").Value = 6684
$ws1.Range("F6").Value = 5481
$ws1.Range("F11").Value = 240
$ws1.Range("F12").Value = 128
$ws1.Range("F13").Value = 45

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 269
$ws4.Range("F5").Value = 6684
$ws4.Range("F6").Value = 5481
$ws4.Range("F11").Value = 240
$ws4.Range("F14").Value = 128
$ws4.Range("F15").Value = 45
